$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zero out the forum view/grade data (columns B:J, rows 2:50) as part of the
# correction of the forum grades for matc65 2021.2.
$ws.Range("B2:J50").Value = 0
